$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value in A8 (e.g., rated energy value) from 1800 to 5000
$ws.Range("A8").Value = 5000

# Update the active cell selection to B10
$ws.Range("B10").Select()
